$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 6467
$ws.Range("J43").Value = 10000
$ws.Range("L43").Value = 10000
$ws.Range("N43").Value = -10138

$ws.Range("H87").Value = 59999.332
$ws.Range("J87").Value = 59999.332
$ws.Range("L87").Value = 59999.332
$ws.Range("N87").Value = -62495.332

$ws.Range("H90").Value = 59999.332
$ws.Range("J90").Value = 59999.332
$ws.Range("L90").Value = 179997.996
$ws.Range("N90").Value = -192477.996

$ws.Range("H98").Value = 31252736
$ws.Range("I98").Value = 33336010
$ws.Range("K98").Value = 33336010
$ws.Range("M98").Value = -33334512

$ws.Range("H121").Value = 5900
$ws.Range("J121").Value = 5900
$ws.Range("L121").Value = 17700
$ws.Range("N121").Value = -21194

$ws.Range("H122").Value = 31252736
$ws.Range("I122").Value = 33336010
$ws.Range("K122").Value = 100008030
$ws.Range("M122").Value = -100005580

$ws.Range("H138").Value = 5405.915
$ws.Range("I138").Value = 1180.9474
$ws.Range("J138").Value = 8272.857
$ws.Range("K138").Value = 3542.8422
$ws.Range("L138").Value = 24818.571
$ws.Range("M138").Value = 1597.1578
$ws.Range("N138").Value = -35098.571

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1840770.2
$ws.Range("I32").Value = 1986502.1
$ws.Range("K32").Value = 1986502.1
$ws.Range("M32").Value = -1986215.1

$ws.Range("H61").Value = 13230.23
$ws.Range("I61").Value = 2999.5
$ws.Range("J61").Value = 17777.223
$ws.Range("K61").Value = 2999.5
$ws.Range("L61").Value = 17777.223
$ws.Range("M61").Value = -2787.5
$ws.Range("N61").Value = -18201.223

$ws.Range("H132").Value = 8425.852999999999
$ws.Range("I132").Value = 6899.1177
$ws.Range("K132").Value = 20697.3531
$ws.Range("M132").Value = -18167.3531

$ws.Range("H136").Value = 13230.23
$ws.Range("I136").Value = 2999.5
$ws.Range("J136").Value = 17777.223
$ws.Range("K136").Value = 8998.5
$ws.Range("L136").Value = 53331.66900000001
$ws.Range("M136").Value = -6448.5
$ws.Range("N136").Value = -58431.66900000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 2000
$ws.Range("I23").Value = 2000
$ws.Range("K23").Value = 2000
$ws.Range("M23").Value = -1760

$ws.Range("H27").Value = 2000
$ws.Range("I27").Value = 2000
$ws.Range("K27").Value = 2000
$ws.Range("M27").Value = -1808

$ws.Range("H59").Value = 97496.5
$ws.Range("J59").Value = 97496.5
$ws.Range("L59").Value = 97496.5
$ws.Range("N59").Value = -99786.5

$ws.Range("H119").Value = 90500
$ws.Range("J119").Value = 90500
$ws.Range("L119").Value = 90500
$ws.Range("N119").Value = -100176

$ws.Range("H122").Value = 1500
$ws.Range("I122").Value = 1500
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4500
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2050
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 4226.706
$ws.Range("I132").Value = 1689.9131
$ws.Range("K132").Value = 5069.7393
$ws.Range("M132").Value = -2539.7393

$ws.Range("H134").Value = 9363.297
$ws.Range("I134").Value = 11395.286
$ws.Range("J134").Value = 8126.4346
$ws.Range("K134").Value = 34185.858
$ws.Range("L134").Value = 24379.3038
$ws.Range("M134").Value = -31650.858
$ws.Range("N134").Value = -29449.3038

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H99").Value = 5991.222
$ws.Range("J99").Value = 19500
$ws.Range("L99").Value = 58500
$ws.Range("N99").Value = -62992

$ws.Range("H108").Value = 6031.75
$ws.Range("I108").Value = 127
$ws.Range("J108").Value = 8000
$ws.Range("K108").Value = 381
$ws.Range("L108").Value = 24000
$ws.Range("M108").Value = 2499
$ws.Range("N108").Value = -29760

$ws.Range("H109").Value = 55556456
$ws.Range("J109").Value = 55557556
$ws.Range("L109").Value = 166672668
$ws.Range("N109").Value = -166674748

$ws.Range("H120").Value = 12115.667
$ws.Range("I120").Value = 6341.3335
$ws.Range("J120").Value = 17890
$ws.Range("K120").Value = 19024.0005
$ws.Range("L120").Value = 53670
$ws.Range("M120").Value = -14186.0005
$ws.Range("N120").Value = -63346

$ws.Range("H131").Value = 1979.8064
$ws.Range("I131").Value = 1315.25
$ws.Range("J131").Value = 2210.9565
$ws.Range("K131").Value = 3945.75
$ws.Range("L131").Value = 6632.869499999999
$ws.Range("M131").Value = 1094.25
$ws.Range("N131").Value = -16712.8695

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 42069
$ws.Range("J26").Value = 42069
$ws.Range("L26").Value = 42069
$ws.Range("N26").Value = -42629

$ws.Range("H50").Value = 42069
$ws.Range("J50").Value = 42069
$ws.Range("L50").Value = 42069
$ws.Range("N50").Value = -43065

$ws.Range("H52").Value = 79989.39999999999
$ws.Range("J52").Value = 79989.39999999999
$ws.Range("L52").Value = 79989.39999999999
$ws.Range("N52").Value = -80507.39999999999

$ws.Range("H58").Value = 59588
$ws.Range("J58").Value = 73724.75
$ws.Range("L58").Value = 73724.75
$ws.Range("N58").Value = -74278.75

$ws.Range("H80").Value = 4395.857
$ws.Range("I80").Value = 3690
$ws.Range("K80").Value = 3690
$ws.Range("M80").Value = -2692

$ws.Range("H83").Value = 4395.857
$ws.Range("I83").Value = 3690
$ws.Range("K83").Value = 18450
$ws.Range("M83").Value = -13458

$ws.Range("H122").Value = 47266.78
$ws.Range("I122").Value = 69689.47
$ws.Range("K122").Value = 209068.41
$ws.Range("M122").Value = -206618.41

$ws.Range("H123").Value = 50001
$ws.Range("J123").Value = 50001
$ws.Range("L123").Value = 50001
$ws.Range("N123").Value = -54901

$ws.Range("H132").Value = 7207.75
$ws.Range("I132").Value = 2943.6667
$ws.Range("K132").Value = 8831.000100000001
$ws.Range("M132").Value = -6301.000100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 16178997
$ws.Range("I46").Value = 11494787
$ws.Range("K46").Value = 11494787
$ws.Range("M46").Value = -11494599

$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()

$ws.Range("H136").Value = 6243.939
$ws.Range("I136").Value = 2616.6
$ws.Range("K136").Value = 7849.799999999999
$ws.Range("M136").Value = -5299.799999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 174999.72
$ws.Range("I15").Value = 174999.72
$ws.Range("K15").Value = 174999.72
$ws.Range("M15").Value = -174711.72

$ws.Range("H132").Value = 11379493
$ws.Range("I132").Value = 14710206
$ws.Range("K132").Value = 44130618
$ws.Range("M132").Value = -44128088

$ws.Range("H136").Value = 55560400
$ws.Range("I136").Value = 111113220
$ws.Range("J136").Value = 7578.6665
$ws.Range("K136").Value = 333339660
$ws.Range("L136").Value = 22735.9995
$ws.Range("M136").Value = -333337110
$ws.Range("N136").Value = -27835.9995
